$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column F. Everything from F
# onward (F..J) shifts right by one (F->G, J->K).
$ws.Columns("F").Insert()

# --- New column F: "Dig Speed Mod (Ring)" ------------------------------
# Header + color-row label (added in this order so the shared-string
# table gets MAGENTA before "Dig Speed Mod (Ring)", matching the source).
$ws.Range("F6").Value = "MAGENTA"
$ws.Range("F1").Value = "Dig Speed Mod (Ring)"

# Hidden helper rows 2-5: same pattern as column E, shifted one column
# right (references B2 -> C2, etc.), mirroring a fill-right from E.
$ws.Range("F2").Formula = "=1+LOG(C2/100,32)"
$ws.Range("F3").Formula = "=2+LOG(C3/100,8)"
$ws.Range("F4").Formula = "=2+LOG(C4/100,8)"
$ws.Range("F5").Formula = "=2+LOG(C5/100,8)"

# Row 7: standalone formula (first row of the visible table).
$ws.Range("F7").Formula = "=1.1 + 0.1 * LOG(B7/100,8)"

# Rows 8-16: shared formula, filled down only through row 16.
$ws.Range("F8:F16").Formula = "=1.1 + 0.1 * LOG(B8/100,8)"

# Highlight F16 (matches the existing highlight pattern used elsewhere
# in the sheet, e.g. C16 / E25).
$ws.Range("F16").Style = "Normal"
$ws.Range("F16").Interior.Color = $ws.Range("C16").Interior.Color

# Column widths: new F column keeps the old F width (now on G); the
# inserted F column gets its own explicit width.
$ws.Columns("F").ColumnWidth = 20.140625
$ws.Columns("G").ColumnWidth = 19.42578125

# Update the active selection to match the saved view.
$ws.Range("F16").Select()
